# Update NATMI LR-pair statistics for Sema4d-Met sheet
# following revised ligand/receptor-expressing cell counts (Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.6496915
$ws.Range("H2").Value = 1.299383
$ws.Range("I2").Value = 0.005425816501278846
$ws.Range("J2").Value = 0.003655611729725819
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 3.180483
$ws.Range("N2").Value = 6.360966
$ws.Range("O2").Value = 0.06092186256078212
$ws.Range("P2").Value = 0.05506293733437207
$ws.Range("Q2").Value = 2.0663327709945
$ws.Range("R2").Value = 8.265331083977999
$ws.Range("S2").Value = 0.0003305508471709335
$ws.Range("T2").Value = 0.0002012887195926883

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.6496915
$ws.Range("H3").Value = 1.299383
$ws.Range("I3").Value = 0.005425816501278846
$ws.Range("J3").Value = 0.003655611729725819
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7202763333333334
$ws.Range("N3").Value = 2.160829
$ws.Range("O3").Value = 0.01379682764696979
$ws.Range("P3").Value = 0.01870495641971579
$ws.Range("Q3").Value = 0.4679574114178334
$ws.Range("R3").Value = 2.807744468507
$ws.Range("S3").Value = 0.00007485905511222888
$ws.Range("T3").Value = 0.00006837805809192331

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 0.6496915
$ws.Range("H4").Value = 1.299383
$ws.Range("I4").Value = 0.005425816501278846
$ws.Range("J4").Value = 0.003655611729725819
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.55178
$ws.Range("N4").Value = 13.65534
$ws.Range("O4").Value = 0.08718893185938011
$ws.Range("P4").Value = 0.1182058087874616
$ws.Range("Q4").Value = 2.95725277587
$ws.Range("R4").Value = 17.74351665522
$ws.Range("S4").Value = 0.0004730711452115014
$ws.Range("T4").Value = 0.0004321145411251718

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 0.6496915
$ws.Range("H5").Value = 1.299383
$ws.Range("I5").Value = 0.005425816501278846
$ws.Range("J5").Value = 0.003655611729725819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.321411
$ws.Range("N5").Value = 12.964233
$ws.Range("O5").Value = 0.08277623461928645
$ws.Range("P5").Value = 0.1122233241408928
$ws.Range("Q5").Value = 2.8075839947065
$ws.Range("R5").Value = 16.845503968239
$ws.Range("S5").Value = 0.0004491286597110537
$ws.Range("T5").Value = 0.0004102449000782705

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 0.6496915
$ws.Range("H6").Value = 1.299383
$ws.Range("I6").Value = 0.005425816501278846
$ws.Range("J6").Value = 0.003655611729725819
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.51639
$ws.Range("N6").Value = 4.54917
$ws.Range("O6").Value = 0.02904631251559728
$ws.Range("P6").Value = 0.03937934311131445
$ws.Range("Q6").Value = 0.9851856936850001
$ws.Range("R6").Value = 5.91111416211
$ws.Range("S6").Value = 0.00015759996174843
$ws.Range("T6").Value = 0.0001439555885866187

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 0.6496915
$ws.Range("H7").Value = 1.299383
$ws.Range("I7").Value = 0.005425816501278846
$ws.Range("J7").Value = 0.003655611729725819
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 37.915598
$ws.Range("N7").Value = 75.83119600000001
$ws.Range("O7").Value = 0.7262698307979842
$ws.Range("P7").Value = 0.6564236302062433
$ws.Range("Q7").Value = 24.633441738017
$ws.Range("R7").Value = 98.53376695206801
$ws.Range("S7").Value = 0.003940606832324698
$ws.Range("T7").Value = 0.002399629922251146

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.9977106666666667
$ws.Range("H8").Value = 2.993132
$ws.Range("I8").Value = 0.008332254614231399
$ws.Range("J8").Value = 0.008420710789519103
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 3.180483
$ws.Range("N8").Value = 6.360966
$ws.Range("O8").Value = 0.06092186256078212
$ws.Range("P8").Value = 0.05506293733437207
$ws.Range("Q8").Value = 3.173201814252
$ws.Range("R8").Value = 19.039210885512
$ws.Range("S8").Value = 0.0005076164704296479
$ws.Range("T8").Value = 0.0004636690705141611

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.9977106666666667
$ws.Range("H9").Value = 2.993132
$ws.Range("I9").Value = 0.008332254614231399
$ws.Range("J9").Value = 0.008420710789519103
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7202763333333334
$ws.Range("N9").Value = 2.160829
$ws.Range("O9").Value = 0.01379682764696979
$ws.Range("P9").Value = 0.01870495641971579
$ws.Range("Q9").Value = 0.7186273807142224
$ws.Range("R9").Value = 6.467646426428001
$ws.Range("S9").Value = 0.0001149586808232194
$ws.Range("T9").Value = 0.0001575090283409854

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.9977106666666667
$ws.Range("H10").Value = 2.993132
$ws.Range("I10").Value = 0.008332254614231399
$ws.Range("J10").Value = 0.008420710789519103
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.55178
$ws.Range("N10").Value = 13.65534
$ws.Range("O10").Value = 0.08718893185938011
$ws.Range("P10").Value = 0.1182058087874616
$ws.Range("Q10").Value = 4.541359458320001
$ws.Range("R10").Value = 40.87223512488001
$ws.Range("S10").Value = 0.0007264803797952269
$ws.Range("T10").Value = 0.0009953769294404097

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.9977106666666667
$ws.Range("H11").Value = 2.993132
$ws.Range("I11").Value = 0.008332254614231399
$ws.Range("J11").Value = 0.008420710789519103
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.321411
$ws.Range("N11").Value = 12.964233
$ws.Range("O11").Value = 0.08277623461928645
$ws.Range("P11").Value = 0.1122233241408928
$ws.Range("Q11").Value = 4.311517849750667
$ws.Range("R11").Value = 38.803660647756
$ws.Range("S11").Value = 0.0006897126628552504
$ws.Range("T11").Value = 0.0009450001564289158

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9977106666666667
$ws.Range("H12").Value = 2.993132
$ws.Range("I12").Value = 0.008332254614231399
$ws.Range("J12").Value = 0.008420710789519103
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.51639
$ws.Range("N12").Value = 4.54917
$ws.Range("O12").Value = 0.02904631251559728
$ws.Range("P12").Value = 0.03937934311131445
$ws.Range("Q12").Value = 1.512918477826667
$ws.Range("R12").Value = 13.61626630044
$ws.Range("S12").Value = 0.0002420212714844927
$ws.Range("T12").Value = 0.0003316020594216203

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9977106666666667
$ws.Range("H13").Value = 2.993132
$ws.Range("I13").Value = 0.008332254614231399
$ws.Range("J13").Value = 0.008420710789519103
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 37.915598
$ws.Range("N13").Value = 75.83119600000001
$ws.Range("O13").Value = 0.7262698307979842
$ws.Range("P13").Value = 0.6564236302062433
$ws.Range("Q13").Value = 37.82879655764534
$ws.Range("R13").Value = 226.972779345872
$ws.Range("S13").Value = 0.006051465148843561
$ws.Range("T13").Value = 0.00552755354537301

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 42.05654266666667
$ws.Range("H14").Value = 126.169628
$ws.Range("I14").Value = 0.3512299040198892
$ws.Range("J14").Value = 0.354958601160661
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 3.180483
$ws.Range("N14").Value = 6.360966
$ws.Range("O14").Value = 0.06092186256078212
$ws.Range("P14").Value = 0.05506293733437207
$ws.Range("Q14").Value = 133.760118990108
$ws.Range("R14").Value = 802.560713940648
$ws.Range("S14").Value = 0.02139757993993638
$ws.Range("T14").Value = 0.01954506321200584

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 42.05654266666667
$ws.Range("H15").Value = 126.169628
$ws.Range("I15").Value = 0.3512299040198892
$ws.Range("J15").Value = 0.354958601160661
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.7202763333333334
$ws.Range("N15").Value = 2.160829
$ws.Range("O15").Value = 0.01379682764696979
$ws.Range("P15").Value = 0.01870495641971579
$ws.Range("Q15").Value = 30.29233234462356
$ws.Range("R15").Value = 272.630991101612
$ws.Range("S15").Value = 0.004845858450224154
$ws.Range("T15").Value = 0.006639485165513443

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 42.05654266666667
$ws.Range("H16").Value = 126.169628
$ws.Range("I16").Value = 0.3512299040198892
$ws.Range("J16").Value = 0.354958601160661
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.55178
$ws.Range("N16").Value = 13.65534
$ws.Range("O16").Value = 0.08718893185938011
$ws.Range("P16").Value = 0.1182058087874616
$ws.Range("Q16").Value = 191.43212977928
$ws.Range("R16").Value = 1722.88916801352
$ws.Range("S16").Value = 0.03062336016856673
$ws.Range("T16").Value = 0.04195816853626193

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 42.05654266666667
$ws.Range("H17").Value = 126.169628
$ws.Range("I17").Value = 0.3512299040198892
$ws.Range("J17").Value = 0.354958601160661
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.321411
$ws.Range("N17").Value = 12.964233
$ws.Range("O17").Value = 0.08277623461928645
$ws.Range("P17").Value = 0.1122233241408928
$ws.Range("Q17").Value = 181.7436061017027
$ws.Range("R17").Value = 1635.692454915324
$ws.Range("S17").Value = 0.02907348894045981
$ws.Range("T17").Value = 0.03983463415465075

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 42.05654266666667
$ws.Range("H18").Value = 126.169628
$ws.Range("I18").Value = 0.3512299040198892
$ws.Range("J18").Value = 0.354958601160661
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 1.51639
$ws.Range("N18").Value = 4.54917
$ws.Range("O18").Value = 0.02904631251559728
$ws.Range("P18").Value = 0.03937934311131445
$ws.Range("Q18").Value = 63.77412073430667
$ws.Range("R18").Value = 573.96708660876
$ws.Range("S18").Value = 0.01020193355698494
$ws.Range("T18").Value = 0.01397803654541789

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 42.05654266666667
$ws.Range("H19").Value = 126.169628
$ws.Range("I19").Value = 0.3512299040198892
$ws.Range("J19").Value = 0.354958601160661
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 37.915598
$ws.Range("N19").Value = 75.83119600000001
$ws.Range("O19").Value = 0.7262698307979842
$ws.Range("P19").Value = 0.6564236302062433
$ws.Range("Q19").Value = 1594.598965019181
$ws.Range("R19").Value = 9567.593790115088
$ws.Range("S19").Value = 0.2550876829637171
$ws.Range("T19").Value = 0.2330032135468111

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 72.21795166666666
$ws.Range("H20").Value = 216.653855
$ws.Range("I20").Value = 0.6031191016683428
$ws.Range("J20").Value = 0.6095218835619034
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 3.180483
$ws.Range("N20").Value = 6.360966
$ws.Range("O20").Value = 0.06092186256078212
$ws.Range("P20").Value = 0.05506293733437207
$ws.Range("Q20").Value = 229.687967570655
$ws.Range("R20").Value = 1378.12780542393
$ws.Range("S20").Value = 0.03674313901962116
$ws.Range("T20").Value = 0.03356206527849752

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 72.21795166666666
$ws.Range("H21").Value = 216.653855
$ws.Range("I21").Value = 0.6031191016683428
$ws.Range("J21").Value = 0.6095218835619034
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.7202763333333334
$ws.Range("N21").Value = 2.160829
$ws.Range("O21").Value = 0.01379682764696979
$ws.Range("P21").Value = 0.01870495641971579
$ws.Range("Q21").Value = 52.01688142731056
$ws.Range("R21").Value = 468.151932845795
$ws.Range("S21").Value = 0.008321130296313376
$ws.Range("T21").Value = 0.01140108026888849

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 72.21795166666666
$ws.Range("H22").Value = 216.653855
$ws.Range("I22").Value = 0.6031191016683428
$ws.Range("J22").Value = 0.6095218835619034
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 4.55178
$ws.Range("N22").Value = 13.65534
$ws.Range("O22").Value = 0.08718893185938011
$ws.Range("P22").Value = 0.1182058087874616
$ws.Range("Q22").Value = 328.7202280373
$ws.Range("R22").Value = 2958.4820523357
$ws.Range("S22").Value = 0.05258531025845169
$ws.Range("T22").Value = 0.07204902722009178

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 72.21795166666666
$ws.Range("H23").Value = 216.653855
$ws.Range("I23").Value = 0.6031191016683428
$ws.Range("J23").Value = 0.6095218835619034
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 4.321411
$ws.Range("N23").Value = 12.964233
$ws.Range("O23").Value = 0.08277623461928645
$ws.Range("P23").Value = 0.1122233241408928
$ws.Range("Q23").Value = 312.0834507298017
$ws.Range("R23").Value = 2808.751056568215
$ws.Range("S23").Value = 0.04992392826307202
$ws.Range("T23").Value = 0.06840257190993503

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 72.21795166666666
$ws.Range("H24").Value = 216.653855
$ws.Range("I24").Value = 0.6031191016683428
$ws.Range("J24").Value = 0.6095218835619034
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 1.51639
$ws.Range("N24").Value = 4.54917
$ws.Range("O24").Value = 0.02904631251559728
$ws.Range("P24").Value = 0.03937934311131445
$ws.Range("Q24").Value = 109.5105797278167
$ws.Range("R24").Value = 985.59521755035
$ws.Range("S24").Value = 0.01751838591118498
$ws.Range("T24").Value = 0.02400257138663885

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 72.21795166666666
$ws.Range("H25").Value = 216.653855
$ws.Range("I25").Value = 0.6031191016683428
$ws.Range("J25").Value = 0.6095218835619034
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 37.915598
$ws.Range("N25").Value = 75.83119600000001
$ws.Range("O25").Value = 0.7262698307979842
$ws.Range("P25").Value = 0.6564236302062433
$ws.Range("Q25").Value = 2738.186823776764
$ws.Range("R25").Value = 16429.12094266058
$ws.Range("S25").Value = 0.4380272079196996
$ws.Range("T25").Value = 0.4001045674978517

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 0.695089
$ws.Range("H26").Value = 2.085267
$ws.Range("I26").Value = 0.005804947988479781
$ws.Range("J26").Value = 0.005866573985353178
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 3.180483
$ws.Range("N26").Value = 6.360966
$ws.Range("O26").Value = 0.06092186256078212
$ws.Range("P26").Value = 0.05506293733437207
$ws.Range("Q26").Value = 2.210718747987
$ws.Range("R26").Value = 13.264312487922
$ws.Range("S26").Value = 0.0003536482435266538
$ws.Range("T26").Value = 0.0003230307957229594

# Row 27
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 0.695089
$ws.Range("H27").Value = 2.085267
$ws.Range("I27").Value = 0.005804947988479781
$ws.Range("J27").Value = 0.005866573985353178
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 0.7202763333333334
$ws.Range("N27").Value = 2.160829
$ws.Range("O27").Value = 0.01379682764696979
$ws.Range("P27").Value = 0.01870495641971579
$ws.Range("Q27").Value = 0.5006561562603333
$ws.Range("R27").Value = 4.505905406343
$ws.Range("S27").Value = 0.00008008986689667952
$ws.Range("T27").Value = 0.0001097340107290696

# Row 28
$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 0.695089
$ws.Range("H28").Value = 2.085267
$ws.Range("I28").Value = 0.005804947988479781
$ws.Range("J28").Value = 0.005866573985353178
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 4.55178
$ws.Range("N28").Value = 13.65534
$ws.Range("O28").Value = 0.08718893185938011
$ws.Range("P28").Value = 0.1182058087874616
$ws.Range("Q28").Value = 3.16389220842
$ws.Range("R28").Value = 28.47502987578
$ws.Range("S28").Value = 0.0005061272146148092
$ws.Range("T28").Value = 0.0006934631227501542

# Row 29
$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 0.695089
$ws.Range("H29").Value = 2.085267
$ws.Range("I29").Value = 0.005804947988479781
$ws.Range("J29").Value = 0.005866573985353178
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 4.321411
$ws.Range("N29").Value = 12.964233
$ws.Range("O29").Value = 0.08277623461928645
$ws.Range("P29").Value = 0.1122233241408928
$ws.Range("Q29").Value = 3.003765250579
$ws.Range("R29").Value = 27.033887255211
$ws.Range("S29").Value = 0.0004805117366471573
$ws.Range("T29").Value = 0.0006583664339548192

# Row 30
$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 0.695089
$ws.Range("H30").Value = 2.085267
$ws.Range("I30").Value = 0.005804947988479781
$ws.Range("J30").Value = 0.005866573985353178
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 1.51639
$ws.Range("N30").Value = 4.54917
$ws.Range("O30").Value = 0.02904631251559728
$ws.Range("P30").Value = 0.03937934311131445
$ws.Range("Q30").Value = 1.05402600871
$ws.Range("R30").Value = 9.48623407839
$ws.Range("S30").Value = 0.0001686123334101715
$ws.Range("T30").Value = 0.0002310218298571342

# Row 31
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 0.695089
$ws.Range("H31").Value = 2.085267
$ws.Range("I31").Value = 0.005804947988479781
$ws.Range("J31").Value = 0.005866573985353178
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 37.915598
$ws.Range("N31").Value = 75.83119600000001
$ws.Range("O31").Value = 0.7262698307979842
$ws.Range("P31").Value = 0.6564236302062433
$ws.Range("Q31").Value = 26.354715098222
$ws.Range("R31").Value = 158.128290589332
$ws.Range("S31").Value = 0.004215958593384309
$ws.Range("T31").Value = 0.003850957792339042

# Row 32
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 3.1237945
$ws.Range("H32").Value = 6.247589
$ws.Range("I32").Value = 0.026087975207778
$ws.Range("J32").Value = 0.01757661877283757
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 3.180483
$ws.Range("N32").Value = 6.360966
$ws.Range("O32").Value = 0.06092186256078212
$ws.Range("P32").Value = 0.05506293733437207
$ws.Range("Q32").Value = 9.935175302743501
$ws.Range("R32").Value = 39.740701210974
$ws.Range("S32").Value = 0.001589328040097342
$ws.Range("T32").Value = 0.0009678202580389027

# Row 33
$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 3.1237945
$ws.Range("H33").Value = 6.247589
$ws.Range("I33").Value = 0.026087975207778
$ws.Range("J33").Value = 0.01757661877283757
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 0.7202763333333334
$ws.Range("N33").Value = 2.160829
$ws.Range("O33").Value = 0.01379682764696979
$ws.Range("P33").Value = 0.01870495641971579
$ws.Range("Q33").Value = 2.249995248546834
$ws.Range("R33").Value = 13.499971491281
$ws.Range("S33").Value = 0.000359931297600134
$ws.Range("T33").Value = 0.0003287698881518851

# Row 34
$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 3.1237945
$ws.Range("H34").Value = 6.247589
$ws.Range("I34").Value = 0.026087975207778
$ws.Range("J34").Value = 0.01757661877283757
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 4.55178
$ws.Range("N34").Value = 13.65534
$ws.Range("O34").Value = 0.08718893185938011
$ws.Range("P34").Value = 0.1182058087874616
$ws.Range("Q34").Value = 14.21882532921
$ws.Range("R34").Value = 85.31295197526001
$ws.Range("S34").Value = 0.002274582692740154
$ws.Range("T34").Value = 0.002077658437792145

# Row 35
$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 3.1237945
$ws.Range("H35").Value = 6.247589
$ws.Range("I35").Value = 0.026087975207778
$ws.Range("J35").Value = 0.01757661877283757
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 4.321411
$ws.Range("N35").Value = 12.964233
$ws.Range("O35").Value = 0.08277623461928645
$ws.Range("P35").Value = 0.1122233241408928
$ws.Range("Q35").Value = 13.4991999140395
$ws.Range("R35").Value = 80.99519948423699
$ws.Range("S35").Value = 0.00215946435654116
$ws.Range("T35").Value = 0.001972506585845052

# Row 36
$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 3.1237945
$ws.Range("H36").Value = 6.247589
$ws.Range("I36").Value = 0.026087975207778
$ws.Range("J36").Value = 0.01757661877283757
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 1.51639
$ws.Range("N36").Value = 4.54917
$ws.Range("O36").Value = 0.02904631251559728
$ws.Range("P36").Value = 0.03937934311131445
$ws.Range("Q36").Value = 4.736890741855
$ws.Range("R36").Value = 28.42134445113
$ws.Range("S36").Value = 0.0007577594807842738
$ws.Range("T36").Value = 0.0006921557013923412

# Row 37
$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 3.1237945
$ws.Range("H37").Value = 6.247589
$ws.Range("I37").Value = 0.026087975207778
$ws.Range("J37").Value = 0.01757661877283757
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 37.915598
$ws.Range("N37").Value = 75.83119600000001
$ws.Range("O37").Value = 0.7262698307979842
$ws.Range("P37").Value = 0.6564236302062433
$ws.Range("Q37").Value = 118.440536496611
$ws.Range("R37").Value = 473.762145986444
$ws.Range("S37").Value = 0.01894690934001493
$ws.Range("T37").Value = 0.01153770790161724

Write-Host "Updated $([int]504) cells across rows 2-37"
